# Update res_bus vm_pu values for Case_4_107 (380 kV case)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.041855452380072
$ws.Range("D2").Value = 1.048112602419637
$ws.Range("E2").Value = 1.045514437125752
$ws.Range("F2").Value = 1.057231649114056
$ws.Range("I2").Value = 1.042796899191724
$ws.Range("J2").Value = 1.046934411610532
$ws.Range("K2").Value = 1.05087335590201
$ws.Range("L2").Value = 1.048282463003418
$ws.Range("M2").Value = 1.059967184412183
$ws.Range("N2").Value = 1.048421178560075
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.042755729240309
$ws.Range("D3").Value = 1.048806046764527
$ws.Range("E3").Value = 1.046363408118827
$ws.Range("F3").Value = 1.058063253164113
$ws.Range("I3").Value = 1.043021081877379
$ws.Range("J3").Value = 1.047481133072909
$ws.Range("K3").Value = 1.051379087190754
$ws.Range("L3").Value = 1.048942792222813
$ws.Range("M3").Value = 1.060612538743194
$ws.Range("N3").Value = 1.048968676429638
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.043338899664495
$ws.Range("D4").Value = 1.049255306037898
$ws.Range("E4").Value = 1.046913708633568
$ws.Range("F4").Value = 1.058602252647805
$ws.Range("I4").Value = 1.043165305437507
$ws.Range("J4").Value = 1.047834880323182
$ws.Range("K4").Value = 1.051706205213374
$ws.Range("L4").Value = 1.049370383415459
$ws.Range("M4").Value = 1.061030370993399
$ws.Range("N4").Value = 1.049322926041582
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.043584214236788
$ws.Range("D5").Value = 1.049444305762767
$ws.Range("E5").Value = 1.047145283040939
$ws.Range("F5").Value = 1.058829060816062
$ws.Range("I5").Value = 1.043225735836936
$ws.Range("J5").Value = 1.047983590176034
$ws.Range("K5").Value = 1.051843694784982
$ws.Range("L5").Value = 1.049550216584783
$ws.Range("M5").Value = 1.06120608486764
$ws.Range("N5").Value = 1.049471847079482
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.043625412373496
$ws.Range("D6").Value = 1.049476047290967
$ws.Range("E6").Value = 1.047184178729506
$ws.Range("F6").Value = 1.05886715533301
$ws.Range("I6").Value = 1.043235870550303
$ws.Range("J6").Value = 1.048008558859009
$ws.Range("K6").Value = 1.051866778046163
$ws.Range("L6").Value = 1.0495804156556
$ws.Range("M6").Value = 1.061235591313063
$ws.Range("N6").Value = 1.049496851220851
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.043342176983188
$ws.Range("D7").Value = 1.049257830948246
$ws.Range("E7").Value = 1.046916802048024
$ws.Range("F7").Value = 1.058605282435721
$ws.Range("I7").Value = 1.043166113702797
$ws.Range("J7").Value = 1.047836867414542
$ws.Range("K7").Value = 1.051708042478429
$ws.Range("L7").Value = 1.049372786066649
$ws.Range("M7").Value = 1.061032718667823
$ws.Range("N7").Value = 1.049324915954839
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.042159573668624
$ws.Range("D8").Value = 1.048346839425196
$ws.Range("E8").Value = 1.045801151188682
$ws.Range("F8").Value = 1.0575125070204
$ws.Range("I8").Value = 1.042872835821707
$ws.Range("J8").Value = 1.047119181737445
$ws.Range("K8").Value = 1.051044295114963
$ws.Range("L8").Value = 1.048505558521197
$ws.Range("M8").Value = 1.060185233379132
$ws.Range("N8").Value = 1.048606211081761
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.04008056918766
$ws.Range("D9").Value = 1.0467458728922
$ws.Range("E9").Value = 1.043842653032802
$ws.Range("F9").Value = 1.05559383644759
$ws.Range("I9").Value = 1.042349656538031
$ws.Range("J9").Value = 1.045854440590317
$ws.Range("K9").Value = 1.049873790165665
$ws.Range("L9").Value = 1.046979857537617
$ws.Range("M9").Value = 1.058693794048116
$ws.Range("N9").Value = 1.047339673857137
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.038697937028925
$ws.Range("D10").Value = 1.045681559676827
$ws.Range("E10").Value = 1.042542074686107
$ws.Range("F10").Value = 1.05431949150734
$ws.Range("I10").Value = 1.041996615704941
$ws.Range("J10").Value = 1.045011290066024
$ws.Range("K10").Value = 1.049092920103481
$ws.Range("L10").Value = 1.045964460890394
$ws.Range("M10").Value = 1.057700888711556
$ws.Range("N10").Value = 1.046495325962388
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.038100058362138
$ws.Range("D11").Value = 1.045221431431283
$ws.Range("E11").Value = 1.04198013698799
$ws.Range("F11").Value = 1.05376883888882
$ws.Range("I11").Value = 1.041842744174246
$ws.Range("J11").Value = 1.044646213467308
$ws.Range("K11").Value = 1.04875468289931
$ws.Range("L11").Value = 1.045525211840571
$ws.Range("M11").Value = 1.057271296622413
$ws.Range("N11").Value = 1.046129730913027
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.037878102654203
$ws.Range("D12").Value = 1.045050630065719
$ws.Range("E12").Value = 1.041771593442602
$ws.Range("F12").Value = 1.053564476085183
$ws.Range("I12").Value = 1.041785439460804
$ws.Range("J12").Value = 1.044510610888785
$ws.Range("K12").Value = 1.048629030387389
$ws.Range("L12").Value = 1.045362120292883
$ws.Range("M12").Value = 1.05711177986975
$ws.Range("N12").Value = 1.045993935763289
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.037925707317364
$ws.Range("D13").Value = 1.045087262500686
$ws.Range("E13").Value = 1.041816318336736
$ws.Range("F13").Value = 1.053608304693938
$ws.Range("I13").Value = 1.041797738292758
$ws.Range("J13").Value = 1.044539697942271
$ws.Range("K13").Value = 1.048655983989945
$ws.Range("L13").Value = 1.045397101015636
$ws.Range("M13").Value = 1.057145994339025
$ws.Range("N13").Value = 1.046023064123727
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.038081708916204
$ws.Range("D14").Value = 1.045207310668457
$ws.Range("E14").Value = 1.041962894927219
$ws.Range("F14").Value = 1.053751942634075
$ws.Range("I14").Value = 1.041838010404498
$ws.Range("J14").Value = 1.044635004445449
$ws.Range("K14").Value = 1.048744296746246
$ws.Range("L14").Value = 1.04551172930641
$ws.Range("M14").Value = 1.057258109830965
$ws.Range("N14").Value = 1.046118505973071
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.038177842995041
$ws.Range("D15").Value = 1.045281291030837
$ws.Range("E15").Value = 1.042053230183889
$ws.Range("F15").Value = 1.053840465831513
$ws.Range("I15").Value = 1.041862803527672
$ws.Range("J15").Value = 1.044693726376152
$ws.Range("K15").Value = 1.048798707044525
$ws.Range("L15").Value = 1.04558236424716
$ws.Range("M15").Value = 1.057327194935138
$ws.Range("N15").Value = 1.046177311295651
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.038737633449875
$ws.Range("D16").Value = 1.045712112315681
$ws.Range("E16").Value = 1.042579394512654
$ws.Range("F16").Value = 1.054356060826628
$ws.Range("I16").Value = 1.042006806584421
$ws.Range("J16").Value = 1.045035519372612
$ws.Range("K16").Value = 1.049115365463098
$ws.Range("L16").Value = 1.045993621459637
$ws.Range("M16").Value = 1.05772940666074
$ws.Range("N16").Value = 1.04651958967737
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.039088993057786
$ws.Range("D17").Value = 1.045982550746275
$ws.Range("E17").Value = 1.042909771637049
$ws.Range("F17").Value = 1.054679788305292
$ws.Range("I17").Value = 1.04209686794633
$ws.Range("J17").Value = 1.045249921419622
$ws.Range("K17").Value = 1.049313966875964
$ws.Range("L17").Value = 1.046251706816963
$ws.Range("M17").Value = 1.057981796144758
$ws.Range("N17").Value = 1.04673429619988
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.039294013253956
$ws.Range("D18").Value = 1.046140362885944
$ws.Range("E18").Value = 1.043102592752129
$ws.Range("F18").Value = 1.054868723632484
$ws.Range("I18").Value = 1.042149302443185
$ws.Range("J18").Value = 1.045374979697011
$ws.Range("K18").Value = 1.049429796515303
$ws.Range("L18").Value = 1.04640228449744
$ws.Range("M18").Value = 1.058129043608883
$ws.Range("N18").Value = 1.046859532074367
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.039363933046267
$ws.Range("D19").Value = 1.046194184567665
$ws.Range("E19").Value = 1.043168359688008
$ws.Range("F19").Value = 1.054933164397186
$ws.Range("I19").Value = 1.042167164817923
$ws.Range("J19").Value = 1.045417621504795
$ws.Range("K19").Value = 1.049469289492909
$ws.Range("L19").Value = 1.046453634484514
$ws.Range("M19").Value = 1.058179256734213
$ws.Range("N19").Value = 1.046902234438409
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.039051287411162
$ws.Range("D20").Value = 1.045953528012181
$ws.Range("E20").Value = 1.042874313109204
$ws.Range("F20").Value = 1.05464504396157
$ws.Range("I20").Value = 1.042087215219466
$ws.Range("J20").Value = 1.045226917977581
$ws.Range("K20").Value = 1.049292659974415
$ws.Range("L20").Value = 1.046224012464763
$ws.Range("M20").Value = 1.057954713723446
$ws.Range("N20").Value = 1.046711260090313
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.038035766926116
$ws.Range("D21").Value = 1.045171956404874
$ws.Range("E21").Value = 1.041919726650744
$ws.Range("F21").Value = 1.053709640020804
$ws.Range("I21").Value = 1.041826155406248
$ws.Range("J21").Value = 1.04460693896019
$ws.Range("K21").Value = 1.048718291277183
$ws.Range("L21").Value = 1.045477972323724
$ws.Range("M21").Value = 1.057225093141892
$ws.Range("N21").Value = 1.046090400631604
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.037397981479261
$ws.Range("D22").Value = 1.044681192892879
$ws.Range("E22").Value = 1.041320612213397
$ws.Range("F22").Value = 1.053122522573875
$ws.Range("I22").Value = 1.0416611493459
$ws.Range("J22").Value = 1.044217152197945
$ws.Range("K22").Value = 1.048357069997241
$ws.Range("L22").Value = 1.045009284633154
$ws.Range("M22").Value = 1.056766658358551
$ws.Range("N22").Value = 1.045700060327448
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.037736015621692
$ws.Range("D23").Value = 1.044941294561613
$ws.Range("E23").Value = 1.04163811195244
$ws.Range("F23").Value = 1.053433668496755
$ws.Range("I23").Value = 1.041748704184573
$ws.Range("J23").Value = 1.044423783382791
$ws.Range("K23").Value = 1.048548568645384
$ws.Range("L23").Value = 1.045257708626928
$ws.Range("M23").Value = 1.057009653769811
$ws.Range("N23").Value = 1.045906984952277
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.03906832472873
$ws.Range("D24").Value = 1.045966641920426
$ws.Range("E24").Value = 1.042890334928079
$ws.Range("F24").Value = 1.054660743093528
$ws.Range("I24").Value = 1.04209157717032
$ws.Range("J24").Value = 1.045237312239319
$ws.Range("K24").Value = 1.049302287680837
$ws.Range("L24").Value = 1.046236526224344
$ws.Range("M24").Value = 1.057966951002252
$ws.Range("N24").Value = 1.046721669113094
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.040617453468873
$ws.Range("D25").Value = 1.047159239711698
$ws.Range("E25").Value = 1.04434808226002
$ws.Range("F25").Value = 1.056089025807752
$ws.Range("I25").Value = 1.042485663703194
$ws.Range("J25").Value = 1.046181409740699
$ws.Range("K25").Value = 1.050176492892016
$ws.Range("L25").Value = 1.047373987470987
$ws.Range("M25").Value = 1.059079128493014
$ws.Range("N25").Value = 1.047667107341215
